# Capstone-Diner.xlsx edit script
# Applies:
#   - D12, D13: 0.75 -> 0.8
#   - Insert a new worklist row (row 19), pushing the "Work slice" / "Responsive & testing"
#     rows down by one, and add a new "Deloy" task row.
#   - Update row heights for the affected rows.
#   - Restore / adjust the B16:B.. and B20/B21:C.. merges around the inserted row.
#   - Update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Simple value tweaks
# ---------------------------------------------------------------------------
$ws.Range("D12").Value2 = 0.8
$ws.Range("D13").Value2 = 0.8

# ---------------------------------------------------------------------------
# 2) Insert a row at 19 (shifts old rows 19..983 down to 20..984).
#    The pre-existing merges (B16:B19 and B20:C20) straddle / sit right below
#    the insertion point, so Excel grows/shifts them automatically to
#    B16:B20 and B21:C21 respectively - no explicit (Un)Merge calls needed.
# ---------------------------------------------------------------------------
$ws.Rows("19:19").Insert()

# ---------------------------------------------------------------------------
# 3) Populate the newly inserted row 19.
#    Columns B & C take on the (pre-edit) row 18 content/format
#    ("Video Youtube" task); columns D & E take on row 17's content/format
#    (the "check" mark slides from column D to column E).
# ---------------------------------------------------------------------------
$ws.Range("B18:C18").Copy()
$ws.Range("B19").PasteSpecial(-4122)

$ws.Range("D17:E17").Copy()
$ws.Range("D19").PasteSpecial(-4122)

$ws.Range("C19").Value2 = $ws.Range("C18").Value2
$ws.Range("E19").Value2 = $ws.Range("E17").Value2

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Edit row 18 in place: new "Deloy" task, and move its "check" mark to
#    column D (reusing the existing column-D "check" style from row 16).
# ---------------------------------------------------------------------------
$ws.Range("C18").Value2 = "Deloy"

$ws.Range("D16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) Row heights
# ---------------------------------------------------------------------------
$ws.Rows("17:17").RowHeight = 19.95
$ws.Rows("19:19").RowHeight = 17.4
$ws.Rows("20:20").RowHeight = 18.6

# ---------------------------------------------------------------------------
# 6) Make sure the trailing empty row exists (row 984) with the same row
#    height as the other trailing blank rows.
# ---------------------------------------------------------------------------
$ws.Rows("984:984").RowHeight = 14.25

# ---------------------------------------------------------------------------
# 7) Selection
# ---------------------------------------------------------------------------
$ws.Range("G14").Select()
